$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 3 to 2
$ws.Range("A2:A11").Value = 2

# Update the active selection to F7
$ws.Range("F7").Select()
